$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the two new rows of MAC-address / machine registration data
$ws.Range("A31").Value = 10001
$ws.Range("B31").Value = 10030
$ws.Range("C31").Value = "eng"
$ws.Range("D31").Value = $true
$ws.Range("E31").Value = "superadmin"
$ws.Range("F31").Value = "now()"

$ws.Range("A32").Value = 10001
$ws.Range("B32").Value = 10031
$ws.Range("C32").Value = "eng"
$ws.Range("D32").Value = $true
$ws.Range("E32").Value = "superadmin"
$ws.Range("F32").Value = "now()"

# Update the selected/active cell as seen in the saved file
$ws.Range("E31").Select()
$excel.ActiveWindow.ScrollRow = 1
